{"js": "// Color-code the TODO bullet list in \"Stevo\" doc:\n//   - \"Opravi\u0165 koncov\u00e9 lom\u00edtko v URL na gride\"                -> amber (FFC000)\n//   - the next four bullets (grid/nav/dropdown/CSV items)      -> green (92D050)\n// The two remaining bullets (Correspondence_view / users_view) stay uncolored.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst AMBER = \"#FFC000\";\nconst GREEN = \"#92D050\";\n\n// Map a recognizable snippet of each bullet's text to the color it should get.\n// (Snippets avoid spots that may use a non-breaking space \\u00A0 in the\n// source document, and comparison also normalizes whitespace just in case.)\nconst colorRules = [\n  { match: \"Opravi\u0165 koncov\u00e9 lom\u00edtko\", color: AMBER },\n  { match: \"Vymyslie\u0165 rie\u0161enie, ako dosta\u0165\", color: GREEN },\n  { match: \"Do navig\u00e1cie prida\u0165 tla\u010d\u00edtko profil\", color: GREEN },\n  { match: \"Dropdown pre filtrovanie eventov\", color: GREEN },\n  { match: \"CSV export, \u010d\u00edsla d\u00e1va\u0165 do\", color: GREEN },\n];\n\nconst normalizeWs = (s) => s.replace(/[\\s\\u00A0]+/g, \" \");\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = normalizeWs(para.text || \"\");\n  for (const rule of colorRules) {\n    if (text.indexOf(normalizeWs(rule.match)) !== -1) {\n      para.font.color = rule.color;\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Color-code the TODO bullet list in \"Stevo\" doc:\n#   - \"Opravi\u0165 koncov\u00e9 lom\u00edtko v URL na gride\"               -> amber (FFC000)\n#   - the next four bullets (grid/nav/dropdown/CSV items)     -> green (92D050)\n# The two remaining bullets (Correspondence_view / users_view) stay uncolored.\n\n$d = $word.ActiveDocument\n\n# Word's Font.Color is a packed BGR integer (VBA RGB()-style: B*65536 + G*256 + R).\n$amber = 49407      # 0xFFC000 -> R=FF,G=C0,B=00 -> 0x00C0FF = 49407\n$green = 5296274    # 0x92D050 -> R=92,G=D0,B=50 -> 0x50D092 = 5296274\n\n# Match on accent-free substrings so the comparison is robust regardless of\n# how the interop layer round-trips Slovak diacritics / non-breaking spaces.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.Contains(\"Opravi\")) {\n        $p.Range.Font.Color = $amber\n    } elseif ($t.Contains(\"Vymysli\") -or $t.Contains(\"Do navig\") -or $t.Contains(\"Dropdown pre filtrovanie\") -or $t.Contains(\"CSV export\")) {\n        $p.Range.Font.Color = $green\n    }\n}\n"}
